# Generalização tabela LateX ok.
# Rename "Página1" -> "base", duplicate it into "sem_acentos" with the
# accented Portuguese text replaced by the unaccented equivalents, and make
# "sem_acentos" the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original (only) sheet to "base" -----------------------
$base = $wb.Worksheets.Item(1)
$base.Name = "base"

# Move the selection on "base" before duplicating it, so the copy inherits
# it (matches the diff: base F23 -> F2).
$base.Range("F2").Select()

# --- 2. Duplicate "base" right after itself, rename to "sem_acentos" -----
$base.Copy([System.Reflection.Missing]::Value, $base)
$semAcentos = $wb.Worksheets.Item(2)
$semAcentos.Name = "sem_acentos"

# --- 3. Strip accents from column A / F text (B..E numbers stay as-is) ---
$semAcentos.Range("A2").Value  = "Desorganizado"
$semAcentos.Range("F2").Value  = "Organizado"

$semAcentos.Range("A3").Value  = "Duvidoso"
$semAcentos.Range("F3").Value  = "Confiavel"

$semAcentos.Range("A4").Value  = "Desagradavel"
$semAcentos.Range("F4").Value  = "Agradavel"

$semAcentos.Range("A5").Value  = "Difícil de encontrar a informacao"
$semAcentos.Range("F5").Value  = "Facil encontrar a informaca "

$semAcentos.Range("A6").Value  = "Difícil navegacao"
$semAcentos.Range("F6").Value  = "Graficos mais modernos, facil navegacao"

$semAcentos.Range("A7").Value  = "Muita navegacao para conseguir as informacoes solicitadas"
$semAcentos.Range("F7").Value  = "Graficos simples"

$semAcentos.Range("A8").Value  = "Informacao detalhada"
$semAcentos.Range("F8").Value  = "Informacao rapida"

$semAcentos.Range("A9").Value  = "Painel lento"
$semAcentos.Range("F9").Value  = "Painel e rapido"

$semAcentos.Range("A10").Value = "Parecem amadores"
$semAcentos.Range("F10").Value = "Profissional"

$semAcentos.Range("A11").Value = "Obscuridade na proveniencia dos dados"
$semAcentos.Range("F11").Value = "Clareza na proveniencia dos dados"

$semAcentos.Range("A12").Value = "Falta de confianca"
$semAcentos.Range("F12").Value = "Confianca"

$semAcentos.Range("A13").Value = "Sem tabela clara de casos/UF"
$semAcentos.Range("F13").Value = "Tabela clara de casos/UF"

# --- 4. Trim the sheet back down to the 13 rows that actually have data --
$semAcentos.Range("A14:F1000").EntireRow.Delete()

# Row 1 on "sem_acentos" goes back to the default row height (no more
# custom 23.85 inherited from "base").
$semAcentos.Rows.Item(1).AutoFit()

# --- 5. Column widths tweak (slightly different from "base") -------------
$semAcentos.Columns.Item(1).ColumnWidth   = 49
$semAcentos.Range("B1:E1").ColumnWidth    = 25
$semAcentos.Columns.Item(6).ColumnWidth   = 34.1666666666667

# --- 6. Page setup (margins + header/footer) matching the diff -----------
$semAcentos.PageSetup.LeftMargin   = 56.7
$semAcentos.PageSetup.RightMargin  = 56.7
$semAcentos.PageSetup.TopMargin    = 75.8
$semAcentos.PageSetup.BottomMargin = 75.8
$semAcentos.PageSetup.HeaderMargin = 56.7
$semAcentos.PageSetup.FooterMargin = 56.7
$semAcentos.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&Kffffff&A'
$semAcentos.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12&KffffffPágina &P'

# --- 7. Make "sem_acentos" the active sheet / selection -------------------
$semAcentos.Activate()
$semAcentos.Range("E18").Select()

# Slight width tweak on "base" too (49.79->49.8, 26.12->26.13 per diff).
$base.Columns.Item(1).ColumnWidth = 49
$base.Range("B1:E1").ColumnWidth  = 25.3333333333333
